$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 201
$ws1.Range("F3").Value = 535
$ws1.Range("F4").Value = 45
$ws1.Range("F7").Value = 35
$ws1.Range("F9").Value = 364
$ws1.Range("F10").Value = 3391
$ws1.Range("F11").Value = 41

$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 94

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 201
$ws4.Range("F3").Value = 94
$ws4.Range("F4").Value = 535
$ws4.Range("F5").Value = 45
$ws4.Range("F8").Value = 35
$ws4.Range("F10").Value = 364
$ws4.Range("F11").Value = 3391
$ws4.Range("F12").Value = 41
